# Update FioranoM.xlsx sheet data: refresh the rolling-7-day sums (columns C/D)
# for rows 90-112, shift the date/new-cases series (columns A/B) down by one
# row starting at row 93, and extend the table with two new trailing rows
# (114 and 115) that mirror the existing "not enough data yet" blank C/D
# pattern used by rows 111-113.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 90-112: updated A (date serial), B (nuovi pos.), C (somma mobile 7gg.),
# D (somma mobile 7gg. per 100mila abitanti). These rows already exist and
# already hold numeric values in C/D, so plain value assignment preserves
# cell formatting/typing.
$data = @(
    @(90,  44232, 1,  12, 70.03209804493727),
    @(91,  44233, 1,  12, 70.03209804493727),
    @(92,  44234, 2,  16, 93.37613072658301),
    @(93,  44235, 2,  16, 93.37613072658301),
    @(94,  44236, 0,  19, 110.8841552378173),
    @(95,  44237, 6,  22, 128.3921797490516),
    @(96,  44238, 4,  24, 140.0641960898745),
    @(97,  44239, 4,  33, 192.5882696235775),
    @(98,  44240, 4,  35, 204.2602859644004),
    @(99,  44241, 4,  31, 180.9162532827546),
    @(100, 44242, 11, 32, 186.752261453166),
    @(101, 44243, 2,  37, 215.9323023052232),
    @(102, 44244, 2,  39, 227.6043186460461),
    @(103, 44245, 5,  47, 274.2923840093376),
    @(104, 44246, 9,  46, 268.4563758389261),
    @(105, 44247, 6,  58, 338.4884738838634),
    @(106, 44248, 12, 61, 355.9964983950978),
    @(107, 44249, 10, 69, 402.6845637583893),
    @(108, 44250, 14, 71, 414.3565800992122),
    @(109, 44251, 5,  86, 501.8967026553837),
    @(110, 44252, 13, 89, 519.404727166618),
    @(111, 44253, 11, 88, 513.5687189962066),
    @(112, 44254, 21, 81, 472.7166618033265)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# Row 113: date/new-cases shift down; C/D stay blank (still not enough data
# for the rolling window) -- simply leave the pre-existing blank C113/D113
# cells untouched so their formatting/emptiness is preserved.
$ws.Cells.Item(113, 1).Value = 44255
$ws.Cells.Item(113, 2).Value = 15

# New trailing rows 114 and 115. These rows do not exist yet, so first
# materialize them (with correct formatting and the blank-C/D pattern) by
# copying the format of the immediately preceding row down onto them, then
# overwrite the A/B values. This keeps the same cell style (s="2") used
# throughout column A without Excel minting a brand-new style entry.
$ws.Range("A113:D113").Copy()
$ws.Range("A114:D114").PasteSpecial(-4122)
$ws.Cells.Item(114, 1).Value = 44256
$ws.Cells.Item(114, 2).Value = 9

$ws.Range("A114:D114").Copy()
$ws.Range("A115:D115").PasteSpecial(-4122)
$ws.Cells.Item(115, 1).Value = 44257
$ws.Cells.Item(115, 2).Value = 7

$excel.CutCopyMode = 0
